$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2831.111
$ws.Range("I64").Value = 2761.111
$ws.Range("J64").Value = 2877.7778
$ws.Range("K64").Value = 2761.111
$ws.Range("L64").Value = 2877.7778
$ws.Range("M64").Value = -2513.111
$ws.Range("N64").Value = -3373.7778
$ws.Range("H67").Value = 2831.111
$ws.Range("I67").Value = 2761.111
$ws.Range("J67").Value = 2877.7778
$ws.Range("K67").Value = 2761.111
$ws.Range("L67").Value = 2877.7778
$ws.Range("M67").Value = -1903.111
$ws.Range("N67").Value = -4593.7778
$ws.Range("H86").Value = 21745682
$ws.Range("I86").Value = 8333.933999999999
$ws.Range("J86").Value = 62503210
$ws.Range("K86").Value = 8333.933999999999
$ws.Range("L86").Value = 62503210
$ws.Range("M86").Value = -7210.933999999999
$ws.Range("N86").Value = -62505456
$ws.Range("H89").Value = 21745682
$ws.Range("I89").Value = 8333.933999999999
$ws.Range("J89").Value = 62503210
$ws.Range("K89").Value = 41669.67
$ws.Range("L89").Value = 312516050
$ws.Range("M89").Value = -36053.67
$ws.Range("N89").Value = -312527282
$ws.Range("H95").Value = 27703
$ws.Range("J95").Value = 27703
$ws.Range("L95").Value = 27703
$ws.Range("N95").Value = -33195
$ws.Range("H100").Value = 62501660
$ws.Range("I100").Value = 1715
$ws.Range("K100").Value = 1715
$ws.Range("M100").Value = -1174
$ws.Range("H103").Value = 10526849
$ws.Range("I103").Value = 406
$ws.Range("J103").Value = 14286293
$ws.Range("K103").Value = 1218
$ws.Range("L103").Value = 42858879
$ws.Range("M103").Value = -632
$ws.Range("N103").Value = -42860051
$ws.Range("H123").Value = 46627.035
$ws.Range("J123").Value = 46627.035
$ws.Range("L123").Value = 46627.035
$ws.Range("N123").Value = -56427.035
$ws.Range("H128").Value = 42658.332
$ws.Range("J128").Value = 42658.332
$ws.Range("L128").Value = 42658.332
$ws.Range("N128").Value = -52618.332
$ws.Range("H129").Value = 857.44446
$ws.Range("I129").Value = 580.5
$ws.Range("J129").Value = 1079
$ws.Range("K129").Value = 1741.5
$ws.Range("L129").Value = 3237
$ws.Range("M129").Value = 3258.5
$ws.Range("N129").Value = -13237
$ws.Range("H134").Value = 43170.582
$ws.Range("J134").Value = 43170.582
$ws.Range("L134").Value = 43170.582
$ws.Range("N134").Value = -53310.582

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3553.4666
$ws.Range("I45").Value = 4185.7144
$ws.Range("J45").Value = 3000.25
$ws.Range("K45").Value = 4185.7144
$ws.Range("L45").Value = 3000.25
$ws.Range("M45").Value = -3808.7144
$ws.Range("N45").Value = -3754.25
$ws.Range("H102").Value = 83334950
$ws.Range("I102").Value = 1333.3334
$ws.Range("K102").Value = 1333.3334
$ws.Range("M102").Value = 288.6666
$ws.Range("H123").Value = 2500000
$ws.Range("J123").Value = 2500000
$ws.Range("L123").Value = 2500000
$ws.Range("N123").Value = -2509800
$ws.Range("H127").Value = 49092
$ws.Range("J127").Value = 49092
$ws.Range("L127").Value = 49092
$ws.Range("N127").Value = -59012
$ws.Range("H133").Value = 29349
$ws.Range("J133").Value = 29349
$ws.Range("L133").Value = 29349
$ws.Range("N133").Value = -34409
$ws.Range("H135").Value = 21399.615
$ws.Range("J135").Value = 21399.615
$ws.Range("L135").Value = 21399.615
$ws.Range("N135").Value = -31539.615

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 40780
$ws.Range("J53").Value = 40780
$ws.Range("L53").Value = 40780
$ws.Range("N53").Value = -41928
$ws.Range("H105").Value = 2031.4783
$ws.Range("I105").Value = 1844
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 1844
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -97
$ws.Range("N105").Value = -7494
$ws.Range("H135").Value = 66885.45
$ws.Range("J135").Value = 66885.45
$ws.Range("L135").Value = 66885.45
$ws.Range("N135").Value = -77025.45

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48728.168
$ws.Range("J20").Value = 48728.168
$ws.Range("L20").Value = 48728.168
$ws.Range("N20").Value = -49200.168
$ws.Range("H30").Value = 48728.168
$ws.Range("J30").Value = 48728.168
$ws.Range("L30").Value = 48728.168
$ws.Range("N30").Value = -48910.168
$ws.Range("H31").Value = 2401.4055
$ws.Range("I31").Value = 1651.7
$ws.Range("J31").Value = 3283.4119
$ws.Range("K31").Value = 1651.7
$ws.Range("L31").Value = 3283.4119
$ws.Range("M31").Value = -1356.7
$ws.Range("N31").Value = -3873.4119
$ws.Range("H34").Value = 2401.4055
$ws.Range("I34").Value = 1651.7
$ws.Range("J34").Value = 3283.4119
$ws.Range("K34").Value = 1651.7
$ws.Range("L34").Value = 3283.4119
$ws.Range("M34").Value = -1449.7
$ws.Range("N34").Value = -3687.4119
$ws.Range("H62").Value = 3595.762
$ws.Range("I62").Value = 3500.3333
$ws.Range("J62").Value = 3834.3333
$ws.Range("K62").Value = 3500.3333
$ws.Range("L62").Value = 3834.3333
$ws.Range("M62").Value = -2876.3333
$ws.Range("N62").Value = -5082.3333
$ws.Range("H65").Value = 3595.762
$ws.Range("I65").Value = 3500.3333
$ws.Range("J65").Value = 3834.3333
$ws.Range("K65").Value = 17501.6665
$ws.Range("L65").Value = 19171.6665
$ws.Range("M65").Value = -14381.6665
$ws.Range("N65").Value = -25411.6665
$ws.Range("H99").Value = 1644.1305
$ws.Range("I99").Value = 1326.4
$ws.Range("K99").Value = 1326.4
$ws.Range("M99").Value = 171.5999999999999
$ws.Range("H109").Value = 10914.286
$ws.Range("J109").Value = 10914.286
$ws.Range("L109").Value = 10914.286
$ws.Range("N109").Value = -12994.286
$ws.Range("H124").Value = 27067.5
$ws.Range("J124").Value = 27067.5
$ws.Range("L124").Value = 27067.5
$ws.Range("N124").Value = -31977.5
$ws.Range("H126").Value = 1644.1305
$ws.Range("I126").Value = 1326.4
$ws.Range("K126").Value = 3979.2
$ws.Range("M126").Value = -1509.2
$ws.Range("H127").Value = 52830
$ws.Range("J127").Value = 51745
$ws.Range("L127").Value = 51745
$ws.Range("N127").Value = -61665
$ws.Range("H128").Value = 48728.168
$ws.Range("J128").Value = 48728.168
$ws.Range("L128").Value = 48728.168
$ws.Range("N128").Value = -58688.168
$ws.Range("H135").Value = 54259.832
$ws.Range("J135").Value = 54259.832
$ws.Range("L135").Value = 54259.832
$ws.Range("N135").Value = -64399.832

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10687.667
$ws.Range("J46").Value = 10687.667
$ws.Range("L46").Value = 10687.667
$ws.Range("N46").Value = -10999.667
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H80").Value = 3229.1
$ws.Range("I80").Value = 2680.8333
$ws.Range("J80").Value = 4051.5
$ws.Range("K80").Value = 2680.8333
$ws.Range("L80").Value = 4051.5
$ws.Range("M80").Value = -1682.8333
$ws.Range("N80").Value = -6047.5
$ws.Range("H83").Value = 3229.1
$ws.Range("I83").Value = 2680.8333
$ws.Range("J83").Value = 4051.5
$ws.Range("K83").Value = 13404.1665
$ws.Range("L83").Value = 20257.5
$ws.Range("M83").Value = -8412.166499999999
$ws.Range("N83").Value = -30241.5
$ws.Range("H124").Value = 49980
$ws.Range("J124").Value = 49980
$ws.Range("L124").Value = 49980
$ws.Range("N124").Value = -59800
$ws.Range("H130").Value = 33742
$ws.Range("J130").Value = 33742
$ws.Range("L130").Value = 33742
$ws.Range("N130").Value = -43782
$ws.Range("H132").Value = 2340.2205
$ws.Range("I132").Value = 1661
$ws.Range("J132").Value = 2913.3125
$ws.Range("K132").Value = 4983
$ws.Range("L132").Value = 8739.9375
$ws.Range("M132").Value = -2453
$ws.Range("N132").Value = -13799.9375
$ws.Range("H133").Value = 17872.5
$ws.Range("J133").Value = 17872.5
$ws.Range("L133").Value = 17872.5
$ws.Range("N133").Value = -27992.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 9792.5
$ws.Range("J96").Value = 9792.5
$ws.Range("L96").Value = 9792.5
$ws.Range("N96").Value = -15284.5
$ws.Range("H100").Value = 1096.5217
$ws.Range("I100").Value = 1096.5217
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1096.5217
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -555.5217
$ws.Range("N100").ClearContents()
$ws.Range("H123").Value = 44439.715
$ws.Range("J123").Value = 44439.715
$ws.Range("L123").Value = 44439.715
$ws.Range("N123").Value = -54239.715
$ws.Range("H128").Value = 53135.57
$ws.Range("J128").Value = 53135.57
$ws.Range("L128").Value = 53135.57
$ws.Range("N128").Value = -63095.57
$ws.Range("H132").Value = 2883.7708
$ws.Range("I132").Value = 2851.8
$ws.Range("J132").Value = 2969.8462
$ws.Range("K132").Value = 8555.400000000001
$ws.Range("L132").Value = 8909.5386
$ws.Range("M132").Value = -6025.400000000001
$ws.Range("N132").Value = -13969.5386

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7611.647
$ws.Range("I81").Value = 11559.7
$ws.Range("J81").Value = 1971.5714
$ws.Range("K81").Value = 23119.4
$ws.Range("L81").Value = 3943.1428
$ws.Range("M81").Value = -22058.4
$ws.Range("N81").Value = -6065.1428
$ws.Range("H84").Value = 7611.647
$ws.Range("I84").Value = 11559.7
$ws.Range("J84").Value = 1971.5714
$ws.Range("K84").Value = 115597
$ws.Range("L84").Value = 19715.714
$ws.Range("M84").Value = -110293
$ws.Range("N84").Value = -30323.714
$ws.Range("H109").Value = 29388.5
$ws.Range("J109").Value = 29388.5
$ws.Range("L109").Value = 29388.5
$ws.Range("N109").Value = -32162.5
$ws.Range("H125").Value = 53892
$ws.Range("J125").Value = 53892
$ws.Range("L125").Value = 53892
$ws.Range("N125").Value = -63732
$ws.Range("H127").Value = 27089.889
$ws.Range("J127").Value = 27089.889
$ws.Range("L127").Value = 27089.889
$ws.Range("N127").Value = -37009.889
